# Replace the inline string "Student" labels in column A (rows 2-6) with
# numeric values, as part of dropping the NaN placeholder student names
# from the grade export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 11
$ws.Range("A3").Value = 12
$ws.Range("A4").Value = 13
$ws.Range("A5").Value = 14
$ws.Range("A6").Value = 15
